$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.406.13'
$ws.Range('E2').Value = '  -1.06%  '
$ws.Range('D3').Value = '1.917.53'
$ws.Range('E3').Value = '  +2.27%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '240.75'
$ws.Range('E5').Value = '  +1.48%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.001'
$ws.Range('E6').Value = '  +0.14%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4700'
$ws.Range('E7').Value = '  -1.65%  '
$ws.Range('B8').Value = 'OKB'
$ws.Range('C8').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '44.41'
$ws.Range('E8').Value = '  +1.60%  '
$ws.Range('B9').Value = 'Cardano'
$ws.Range('C9').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2858'
$ws.Range('E9').Value = '  +0.85%  '
$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06988'
$ws.Range('E10').Value = '  +7.31%  '
$ws.Range('B11').Value = 'Litecoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '106.50'
$ws.Range('E11').Value = '  +12.76%  '
$ws.Range('B12').Value = 'Solana'
$ws.Range('C12').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '18.24'
$ws.Range('E12').Value = '  -3.23%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.917.94'
$ws.Range('E13').Value = '  +2.32%  '
$ws.Range('B14').Value = 'TRON'
$ws.Range('C14').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.07650'
$ws.Range('E14').Value = '  +1.86%  '
$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.186'
$ws.Range('E15').Value = '  +1.61%  '
$ws.Range('B16').Value = 'Polygon'
$ws.Range('C16').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.6579'
$ws.Range('E16').Value = '  +1.11%  '
$ws.Range('B17').Value = 'BitcoinCash'
$ws.Range('C17').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '290.87'
$ws.Range('E17').Value = '  -1.30%  '
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').Value = '30.412.80'
$ws.Range('E18').Value = '  -0.96%  '
$ws.Range('B19').Value = 'ShibaInu'
$ws.Range('C19').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007628'
$ws.Range('E19').Value = '  +1.12%  '
$ws.Range('B20').Value = 'Avalanche'
$ws.Range('C20').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.01'
$ws.Range('E20').Value = '  -0.42%  '
$ws.Range('B21').Value = 'Dai'
$ws.Range('C21').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.000'
$ws.Range('E21').Value = '  +0.08%  '
$ws.Range('B22').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C22').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D22').Value = '2.166.84'
$ws.Range('E22').Value = '  +2.67%  '
$ws.Range('B23').Value = 'BinanceUSD'
$ws.Range('C23').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.002'
$ws.Range('E23').Value = '  +0.08%  '
$ws.Range('B24').Value = 'Uniswap'
$ws.Range('C24').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.242'
$ws.Range('E24').Value = '  +1.27%  '
$ws.Range('B25').Value = 'Chainlink'
$ws.Range('C25').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '6.190'
$ws.Range('E25').Value = '  +1.50%  '
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '168.30'
$ws.Range('E26').Value = '  -0.51%  '
$ws.Range('B27').Value = 'Cosmos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.264'
$ws.Range('E27').Value = '  +0.30%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '21.45'
$ws.Range('E28').Value = '  +8.41%  '
$ws.Range('B29').Value = 'LidoDAOToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.047'
$ws.Range('E29').Value = '  +3.43%  '
$ws.Range('B30').Value = 'Stellar'
$ws.Range('C30').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.1080'
$ws.Range('E30').Value = '  +2.62%  '
$ws.Range('B31').Value = 'Toncoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.369'
$ws.Range('E31').Value = '  +1.26%  '
$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.140'
$ws.Range('E32').Value = '  +0.26%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.975'
$ws.Range('E33').Value = '  +0.75%  '
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05068'
$ws.Range('E34').Value = '  +1.86%  '
$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7403'
$ws.Range('E35').Value = '  +2.34%  '
$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.144'
$ws.Range('E36').Value = '  -2.78%  '
$ws.Range('B37').Value = 'HuobiToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.732'
$ws.Range('E37').Value = '  +1.12%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02021'
$ws.Range('E38').Value = '  +4.85%  '
$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.682'
$ws.Range('E39').Value = '  -0.85%  '
$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.065'
$ws.Range('E40').Value = '  +0.52%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8740'
$ws.Range('E41').Value = '  -1.74%  '
$ws.Range('B42').Value = 'Quant'
$ws.Range('C42').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '108.46'
$ws.Range('E42').Value = '  +0.80%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.847'
$ws.Range('E43').Value = '  +4.76%  '
$ws.Range('B44').Value = 'BitcoinSV'
$ws.Range('C44').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '53.69'
$ws.Range('E44').Value = '  +28.23%  '
$ws.Range('B45').Value = 'PaxDollar'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.001'
$ws.Range('E45').Value = '  +0.12%  '
$ws.Range('B46').Value = 'TheSandbox'
$ws.Range('C46').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4216'
$ws.Range('E46').Value = '  +0.25%  '
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '67.53'
$ws.Range('E47').Value = '  +3.24%  '
$ws.Range('B48').Value = 'Aptos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.187'
$ws.Range('E48').Value = '  -2.16%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.242'
$ws.Range('E49').Value = '  +5.33%  '
$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.1207'
$ws.Range('E50').Value = '  -1.73%  '
$ws.Range('B51').Value = 'Elrond'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '34.66'
$ws.Range('E51').Value = '  +0.16%  '
